$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the last existing data row (A205, style with border/bold/centered)
# down onto the new rows A206:A217 so the appended rows match formatting of prior ones.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)

# New "X" index values (column A) for the appended year data (months 204..215 -> 1..12 of next year)
$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)

# Corresponding normalized value data (column B)
$bValues = @(
    "3.469446951953614E-17",
    "7.065055611250996E-17",
    "-1.165734175856414E-16",
    "-1.541976423090495E-18",
    "-2.255140518769849E-17",
    "-8.425799740458777E-17",
    "3.353798720221827E-17",
    "0",
    "5.204170427930421E-17",
    "1.156482317317871E-18",
    "-1.0321604682062E-16",
    "0"
)

$startRow = 206
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = [double]$bValues[$i]
}
